# Metric_Steel_Pan_Head_Screws_with_Split_Lock_Washer.xlsx
#
# The sheet gains two new rows at the top:
#   - new row 1: a plain 0-based numeric column index (0,1,2,...,13)
#     that takes over the bold / bordered "header" styling
#   - new row 2: a mostly blank row with "Washer" in column E
# The old text header row (Lg.,mm, Threading, ...) slides down to row 3,
# loses its bold styling, and its M/N captions (thread_size /
# material_surface) are cleared. All the data rows below it slide down
# by two rows, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 1. This pushes the header
# row down to row 3 (carrying its existing bold/border style with it)
# and the data rows down to rows 4-14.
$ws.Range("A1:A2").EntireRow.Insert()

# Row 1 should be the one styled like a header from now on, so move the
# header formatting (currently sitting on row 3, A:N only) up to row 1.
$ws.Range("A3:N3").Copy()
$ws.Range("A1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The relocated header-text row (row 3) goes back to plain/default
# formatting.
$ws.Rows.Item(3).ClearFormats()

# That row also loses its M/N captions (thread_size / material_surface).
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""

# New row 1 becomes a plain 0-based numeric index, one per column.
for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# New row 2 is blank except for "Washer" in column E.
$ws.Range("E2").Value = "Washer"
